$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (new tpm calculation)

# Row 2
$ws.Range("G2").Value = 0.01135533333333333
$ws.Range("H2").Value = 0.034066
$ws.Range("M2").Value = 13.929953
$ws.Range("N2").Value = 41.789859
$ws.Range("O2").Value = 0.09674275490334808
$ws.Range("P2").Value = 0.09674275490334808
$ws.Range("Q2").Value = 0.1581792596326667
$ws.Range("R2").Value = 1.423613336694
$ws.Range("S2").Value = 0.09674275490334808
$ws.Range("T2").Value = 0.09674275490334808

# Row 3
$ws.Range("G3").Value = 0.01135533333333333
$ws.Range("H3").Value = 0.034066
$ws.Range("M3").Value = 81.07766966666667
$ws.Range("N3").Value = 243.233009
$ws.Range("O3").Value = 0.5630799418129374
$ws.Range("P3").Value = 0.5630799418129373
$ws.Range("Q3").Value = 0.9206639649548889
$ws.Range("R3").Value = 8.285975684594
$ws.Range("S3").Value = 0.5630799418129374
$ws.Range("T3").Value = 0.5630799418129373

# Row 4
$ws.Range("G4").Value = 0.01135533333333333
$ws.Range("H4").Value = 0.034066
$ws.Range("M4").Value = 48.98200233333333
$ws.Range("N4").Value = 146.946007
$ws.Range("O4").Value = 0.3401773032837146
$ws.Range("P4").Value = 0.3401773032837146
$ws.Range("Q4").Value = 0.5562069638291112
$ws.Range("R4").Value = 5.005862674462001
$ws.Range("S4").Value = 0.3401773032837146
$ws.Range("T4").Value = 0.3401773032837146
